$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add boolean column G: row1 = TRUE, rows 2-7 = FALSE
$ws.Range("G1").Value = $true
$ws.Range("G2").Value = $false
$ws.Range("G3").Value = $false
$ws.Range("G4").Value = $false
$ws.Range("G5").Value = $false
$ws.Range("G6").Value = $false
$ws.Range("G7").Value = $false

# Update the selection to match the new active range
$ws.Range("G1:G7").Select()
